$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at M (13): pushes old M (prop_spread) to N,
#    and creates a blank M column for the new Monthly_salary data.
$ws.Columns.Item(13).Insert()

# 2. Header for the new column.
$ws.Cells.Item(1, 13).Value = "Monthly_salary"

# 3. Monthly_salary values for rows 2..358 (municipality monthly salary, NOK).
$monthlySalary = @(65900,57360,68170,59760,62030,55420,52850,55250,54250,57680,59110,56490,67920,63280,58520,55800,54500,56550,59460,59670,58920,58220,63170,55890,57120,57800,58100,56310,55350,57490,59950,55420,54120,51840,53650,55580,57240,54160,53830,55790,57720,53590,53780,58000,53570,54860,57710,55450,51530,54120,54770,58440,55930,52450,53310,54840,53350,54710,54650,55170,53140,54870,50170,49190,55160,53670,53740,55660,57180,50870,55720,55570,53430,51370,54350,54740,53790,55810,52330,53360,50230,50970,51570,53040,52810,54300,52530,52470,54470,54180,48910,55490,55470,58280,54260,56710,59290,55790,55650,53250,53520,52480,53220,52180,76170,69150,59910,65160,56350,65170,65260,60460,60530,55610,60400,60090,53650,54820,61000,62260,56180,53850,54330,55310,54120,57890,61950,54350,63090,62170,55870,53460,52680,51400,52460,51770,54330,53180,52120,51990,54310,50890,51010,54520,59110,58200,55960,53530,51980,55300,51920,53430,51240,50940,51100,49650,55300,50440,55570,50260,49220,48880,51730,51800,51420,51440,51120,51180,51430,50260,51350,52070,54390,51780,52640,50730,52960,51680,53310,53310,54540,52170,50900,51360,50650,51220,50740,50900,50410,56870,57150,58410,56530,55710,61360,59070,56630,54600,55150,58110,55020,54420,52880,53770,55330,53150,54870,56320,53640,53640,56590,55150,55260,56000,57440,57550,55570,55290,56360,52250,54230,54760,53590,57680,52840,53620,51880,52280,52990,54960,52300,52880,54450,53180,54570,54590,54920,61420,56040,55250,57020,58960,58610,56370,59490,56890,55340,55840,53470,54940,56990,56840,58560,59370,58580,58770,54950,53640,53920,57710,58850,61200,55350,54720,54650,53530,55790,53770,54650,52480,53940,58990,52630,54690,55100,55730,54370,53940,53840,51000,60880,54010,55400,54870,51780,53330,50190,52180,52490,52010,55570,56450,61200,53460,55470,51080,57150,51120,55620,53770,51840,51970,52170,51960,52880,52770,54190,53610,53550,54800,52970,55360,53940,56410,53690,53810,55510,52170,58230,57230,53710,55050,54030,50900,51250,56320,53810,54860,54200,52450,53670,52580,50870,52340,51460,53390,52640,53880,53350,55570,57600,55160,55640,53120,52500,50540,49500,50620,52840,54400,52140,49620,52700,50840,49730,52250,52820)
for ($i = 0; $i -lt $monthlySalary.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $monthlySalary[$i]
}

# 4. Fix mangled/inconsistent Region_Name text (column C) so that all
#    municipalities belonging to the same region share one canonical name.
for ($r = 26; $r -le 52; $r++) {
    $ws.Cells.Item($r, 3).Value = "Møre og Romsdal"
}
for ($r = 94; $r -le 105; $r++) {
    $ws.Cells.Item($r, 3).Value = "Østfold"
}
for ($r = 282; $r -le 319; $r++) {
    $ws.Cells.Item($r, 3).Value = "Trøndelag"
}

# 5. Backfill the Grensehandel (column K) constant for rows of the same
#    region that were previously missing it.
$moreExclude = @(30,33,38,41,43,44,48,50)
for ($r = 26; $r -le 52; $r++) {
    if ($moreExclude -notcontains $r) {
        $ws.Cells.Item($r, 11).Value = 8.921212499999999
    }
}
for ($r = 94; $r -le 105; $r++) {
    $ws.Cells.Item($r, 11).Value = 203.3475
}
for ($r = 282; $r -le 319; $r++) {
    $ws.Cells.Item($r, 11).Value = 58.06874999999999
}
